$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the title of the last statistics block: it was wrongly reusing
# "Statistiques sur ARTICLES" (shared string 0); it should read
# "Statistiques sur Sub Ligne".
$ws.Range("A66").Value = "Statistiques sur Sub Ligne"

# The "Sub Ligne" block's header row and first two data rows were left
# empty; fill them in with the column headers and the corresponding
# "Unique values" / "Taux de valeurs manquantes" figures, matching the
# layout already used by the other statistic blocks (e.g. row 42-44).
$ws.Range("B71").Value = "X.ff..fe.IDTICKET"
$ws.Range("C71").Value = "NUMLIGNETICKET"
$ws.Range("D71").Value = "IDARTICLE"

# Copy the header formatting (bold, wrap text) from an equivalent,
# already-styled header row instead of re-creating the style, so it
# reuses the existing style/font table entries.
$ws.Range("B42:D42").Copy()
$ws.Range("B71:D71").PasteSpecial(-4122)

$ws.Range("B72").Value = 326361.0
$ws.Range("C72").Value = 92.0
$ws.Range("D72").Value = 49274.0

$ws.Range("B73").Value = 0.0
$ws.Range("C73").Value = 0.0
$ws.Range("D73").Value = 0.0
